$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "[0. 0. 1.]"
$ws.Range("C3").Value = "[1. 0. 0.]"
$ws.Range("D3").Value = -2
